# Add a new header row at the top of Sheet1 and fill it with the
# column titles used by the app that consumes this sheet
# (fulln, nickn, phonen, officen, listImage, id). All the existing
# data rows shift down by one row (row 2 .. row 6).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a brand-new row 1; everything currently in rows 1-5 moves to rows 2-6
$ws.Rows.Item(1).Insert()

# Fill in the new header row. "id" is written before "listImage" so the
# shared-string table picks up the same ordering as the source edit.
$ws.Range("A1").Value = "fulln"
$ws.Range("B1").Value = "nickn"
$ws.Range("C1").Value = "phonen"
$ws.Range("D1").Value = "officen"
$ws.Range("F1").Value = "id"
$ws.Range("E1").Value = "listImage"

# The new "listImage" column holds long URLs; best-fit the column width
# the same way Excel would after this kind of edit.
$ws.Columns.Item(5).AutoFit()

Write-Output "Inserted header row and re-fit column E"
